$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 132
$ws.Range("H132").Value = 22361.936
$ws.Range("I132").Value = 2987
$ws.Range("J132").Value = 130308
$ws.Range("K132").Value = 8961
$ws.Range("L132").Value = 390924
$ws.Range("M132").Value = -6431
$ws.Range("N132").Value = -395984

# Row 135
$ws.Range("H135").Value = 35715690
$ws.Range("I135").Value = 1639.8889
$ws.Range("J135").Value = 100000980
$ws.Range("K135").Value = 14759.0001
$ws.Range("L135").Value = 900008820
$ws.Range("M135").Value = -12224.0001
$ws.Range("N135").Value = -900013890

# Row 138
$ws.Range("H138").Value = 2332.9
$ws.Range("I138").Value = 1332.75
$ws.Range("J138").Value = 2999.6667
$ws.Range("K138").Value = 3998.25
$ws.Range("L138").Value = 8999.000100000001
$ws.Range("M138").Value = 1141.75
$ws.Range("N138").Value = -19279.0001

# Row 141
$ws.Range("H141").Value = 3572.3076
$ws.Range("I141").Value = 1947.1428
$ws.Range("J141").Value = 5468.3335
$ws.Range("K141").Value = 5841.428400000001
$ws.Range("L141").Value = 16405.0005
$ws.Range("M141").Value = -661.4284000000007
$ws.Range("N141").Value = -26765.0005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 34
$ws.Range("H34").Value = 23000
$ws.Range("J34").Value = 23000
$ws.Range("L34").Value = 23000
$ws.Range("N34").Value = -23542

# Row 74
$ws.Range("H74").Value = 1488.0682
$ws.Range("I74").Value = 1216.0952
$ws.Range("K74").Value = 1216.0952
$ws.Range("M74").Value = -342.0952

# Row 77
$ws.Range("H77").Value = 1488.0682
$ws.Range("I77").Value = 1216.0952
$ws.Range("K77").Value = 6080.476
$ws.Range("M77").Value = -1712.476

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 115
$ws.Range("H115").Value = 20000
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 20000
$ws.Range("N115").Value = -23134

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# Row 129
$ws.Range("H129").Value = 450000
$ws.Range("J129").Value = 450000
$ws.Range("L129").Value = 450000
$ws.Range("N129").Value = -460000

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 9435159
$ws.Range("I58").Value = 729.4857
$ws.Range("J58").Value = 27779884
$ws.Range("K58").Value = 729.4857
$ws.Range("L58").Value = 27779884
$ws.Range("M58").Value = -526.4857
$ws.Range("N58").Value = -27780290

# Row 107
$ws.Range("H107").Value = 395.9
$ws.Range("I107").Value = 395.9
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 395.9
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1524.1
$ws.Range("N107").ClearContents()

# Row 136
$ws.Range("H136").Value = 9435159
$ws.Range("I136").Value = 729.4857
$ws.Range("J136").Value = 27779884
$ws.Range("K136").Value = 2188.4571
$ws.Range("L136").Value = 83339652
$ws.Range("M136").Value = 361.5429000000004
$ws.Range("N136").Value = -83344752

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 994.6391599999999
$ws.Range("I131").Value = 859
$ws.Range("J131").Value = 998.9681
$ws.Range("K131").Value = 2577
$ws.Range("L131").Value = 2996.9043
$ws.Range("M131").Value = 2463
$ws.Range("N131").Value = -13076.9043

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 24
$ws.Range("H24").Value = 37502310
$ws.Range("I24").Value = 100000800
$ws.Range("J24").Value = 3220
$ws.Range("K24").Value = 100000800
$ws.Range("L24").Value = 3220
$ws.Range("M24").Value = -100000627
$ws.Range("N24").Value = -3566

# Row 132
$ws.Range("H132").Value = 3044.5
$ws.Range("I132").Value = 2370
$ws.Range("K132").Value = 7110
$ws.Range("M132").Value = -4580

# Row 140
$ws.Range("H140").Value = 37640
$ws.Range("J140").Value = 37640
$ws.Range("L140").Value = 37640
$ws.Range("N140").Value = -48000

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 71431704
$ws.Range("I7").Value = 100002344
$ws.Range("J7").Value = 5121.25
$ws.Range("K7").Value = 100002344
$ws.Range("L7").Value = 5121.25
$ws.Range("M7").Value = -100002232
$ws.Range("N7").Value = -5345.25

# Row 29
$ws.Range("H29").Value = 26000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Row 40
$ws.Range("H40").Value = 4798.696
$ws.Range("I40").Value = 4386.353
$ws.Range("K40").Value = 4386.353
$ws.Range("M40").Value = -4250.353

# Row 93
$ws.Range("H93").Value = 913.9524
$ws.Range("I93").Value = 466
$ws.Range("J93").Value = 1249.9166
$ws.Range("K93").Value = 466
$ws.Range("L93").Value = 1249.9166
$ws.Range("M93").Value = 782
$ws.Range("N93").Value = -3745.9166

# Row 122
$ws.Range("H122").Value = 2252.2222
$ws.Range("I122").Value = 2252.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6756.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4306.6666
$ws.Range("N122").ClearContents()

# Row 126
$ws.Range("H126").Value = 71431704
$ws.Range("I126").Value = 100002344
$ws.Range("J126").Value = 5121.25
$ws.Range("K126").Value = 300007032
$ws.Range("L126").Value = 15363.75
$ws.Range("M126").Value = -300004562
$ws.Range("N126").Value = -20303.75

# Row 132
$ws.Range("H132").Value = 2684.8718
$ws.Range("I132").Value = 2096.6428
$ws.Range("J132").Value = 4182.1816
$ws.Range("K132").Value = 6289.928400000001
$ws.Range("L132").Value = 12546.5448
$ws.Range("M132").Value = -3759.928400000001
$ws.Range("N132").Value = -17606.5448

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 29
$ws.Range("H29").Value = 41319.418
$ws.Range("I29").Value = 3500
$ws.Range("J29").Value = 48883.3
$ws.Range("K29").Value = 3500
$ws.Range("L29").Value = 48883.3
$ws.Range("M29").Value = -3210
$ws.Range("N29").Value = -49463.3

# Row 32
$ws.Range("H32").Value = 29000
$ws.Range("J32").Value = 29000
$ws.Range("L32").Value = 29000
$ws.Range("N32").Value = -29634

# Row 34
$ws.Range("H34").Value = 17000
$ws.Range("J34").Value = 17000
$ws.Range("L34").Value = 17000
$ws.Range("N34").Value = -17406

# Row 46
$ws.Range("H46").Value = 102262.5
$ws.Range("J46").Value = 102262.5
$ws.Range("L46").Value = 102262.5
$ws.Range("N46").Value = -102724.5

# Row 134
$ws.Range("H134").Value = 102262.5
$ws.Range("J134").Value = 102262.5
$ws.Range("L134").Value = 306787.5
$ws.Range("N134").Value = -311857.5

# Row 136
$ws.Range("H136").Value = 209238.02
$ws.Range("I136").Value = 250710.55
$ws.Range("J136").Value = 1875.375
$ws.Range("K136").Value = 752131.6499999999
$ws.Range("L136").Value = 5626.125
$ws.Range("M136").Value = -749581.6499999999
$ws.Range("N136").Value = -10726.125

# Row 140
$ws.Range("H140").Value = 33931.92
$ws.Range("J140").Value = 33931.92
$ws.Range("L140").Value = 33931.92
$ws.Range("N140").Value = -44291.92

